$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtering save games) - update TB, d2S, K, IP, sum columns
$ws.Range("B2").Value = 0.2881169905109251
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("G2").Value = 2818.197719814496

$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.15379541431027

$ws.Range("B4").Value = 1.445647641019636
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 13.86384647080068
$ws.Range("G4").Value = 17.65757632934944

$ws.Range("B5").Value = 1.445647641019636
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 3.223369029078222
$ws.Range("E5").Value = 13.86384647080068
$ws.Range("G5").Value = 20.15985084044064

$ws.Range("B6").Value = 0.2881169905109251
$ws.Range("C6").Value = 0.3048912486333797
$ws.Range("D6").Value = 0.1496068669990043
$ws.Range("E6").Value = 13.86384647080068
$ws.Range("G6").Value = 14.60646157694399

$ws.Range("B7").Value = 0.1169995834814548
$ws.Range("C7").Value = 1.626987699542094
$ws.Range("D7").Value = 0.1496068669990043
$ws.Range("E7").Value = 2797.565817734744
$ws.Range("G7").Value = 2799.459411884766

$ws.Range("B8").Value = 1.445647641019636
$ws.Range("C8").Value = 0.04103571897497393
$ws.Range("D8").Value = 0.1496068669990043
$ws.Range("E8").Value = 0.5333859586016987
$ws.Range("G8").Value = 2.169676185595313

$ws.Range("B9").Value = 3.272327238179451
$ws.Range("C9").Value = 2919.202174992006
$ws.Range("D9").Value = 3.223369029078222
$ws.Range("E9").Value = 13.86384647080068
$ws.Range("G9").Value = 2939.561717730064
